$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so values like "1.001" or
# percentage strings are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.151.49'
$ws.Range("E2").Value = '  -3.17%  '
$ws.Range("D3").Value = '1.914.04'
$ws.Range("E3").Value = '  -4.04%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.19%  '
$ws.Range("D5").Value = '328.42'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").Value = '0.4666'
$ws.Range("E7").Value = '  -6.05%  '
$ws.Range("D8").Value = '0.4012'
$ws.Range("E8").Value = '  -4.20%  '
$ws.Range("D9").Value = '53.07'
$ws.Range("E9").Value = '  -3.69%  '
$ws.Range("D10").Value = '0.08406'
$ws.Range("E10").Value = '  -5.46%  '
$ws.Range("D11").Value = '1.045'
$ws.Range("E11").Value = '  -4.20%  '
$ws.Range("D12").Value = '22.13'
$ws.Range("E12").Value = '  -3.42%  '
$ws.Range("D13").Value = '1.914.91'
$ws.Range("E13").Value = '  -4.54%  '
$ws.Range("D14").Value = '7.428'
$ws.Range("E14").Value = '  -7.06%  '
$ws.Range("D15").Value = '6.059'
$ws.Range("E15").Value = '  -5.51%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '89.57'
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("D18").Value = '0.00001063'
$ws.Range("E18").Value = '  -3.77%  '
$ws.Range("D19").Value = '0.06627'
$ws.Range("E19").Value = '  -1.65%  '
$ws.Range("D20").Value = '17.92'
$ws.Range("E20").Value = '  -7.82%  '
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = '5.745'
$ws.Range("E22").Value = '  -3.76%  '
$ws.Range("D23").Value = '28.143.52'
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").Value = '11.18'
$ws.Range("E24").Value = '  -6.53%  '
$ws.Range("D25").Value = '2.299'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = '2.140.25'
$ws.Range("E26").Value = '  -4.63%  '
$ws.Range("D27").Value = '153.28'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '20.03'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '5.769'
$ws.Range("E29").Value = '  -7.86%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '2.137'
$ws.Range("E30").Value = '  -4.85%  '
$ws.Range("D31").Value = '123.37'
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").Value = '0.9770'
$ws.Range("E32").Value = '  -6.23%  '
$ws.Range("D33").Value = '0.09655'
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").Value = '1.439'
$ws.Range("E34").Value = '  -5.99%  '
$ws.Range("D35").Value = '3.648'
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("E36").Value = '  -4.85%  '
$ws.Range("D37").Value = '8.868'
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '1.268'
$ws.Range("E38").Value = '  -3.23%  '
$ws.Range("D39").Value = '0.02297'
$ws.Range("E39").Value = '  -4.90%  '
$ws.Range("D40").Value = '0.06177'
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D41").Value = '0.6164'
$ws.Range("E41").Value = '  -4.61%  '
$ws.Range("D42").Value = '11.04'
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").Value = '0.1903'
$ws.Range("E44").Value = '  -3.58%  '
$ws.Range("D45").Value = '1.306'
$ws.Range("E45").Value = '  -3.56%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '12.79'
$ws.Range("E46").Value = '  -4.43%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5858'
$ws.Range("E47").Value = '  -5.15%  '
$ws.Range("D48").Value = '2.023'
$ws.Range("E48").Value = '  -6.59%  '
$ws.Range("D49").Value = '3.436'
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("D50").Value = '0.06907'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '111.56'
$ws.Range("E51").Value = '  -1.06%  '
